$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 53.333332
$ws.Range("I6").Value = 53.333332
$ws.Range("K6").Value = 159.999996
$ws.Range("M6").Value = -47.99999600000001

$ws.Range("H19").Value = 887.7778
$ws.Range("I19").Value = 818
$ws.Range("K19").Value = 818
$ws.Range("M19").Value = -643

$ws.Range("H40").Value = 1991.25
$ws.Range("I40").Value = 1320
$ws.Range("J40").Value = 2662.5
$ws.Range("K40").Value = 1320
$ws.Range("L40").Value = 2662.5
$ws.Range("M40").Value = -1145
$ws.Range("N40").Value = -3012.5

$ws.Range("H74").Value = 11368245
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 12504720
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 12504720
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -12506592

$ws.Range("H77").Value = 11368245
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 12504720
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 62523600
$ws.Range("M77").Value = -12820
$ws.Range("N77").Value = -62532960

$ws.Range("H82").Value = 2042
$ws.Range("I82").Value = 2042
$ws.Range("K82").Value = 6126
$ws.Range("M82").Value = -5720

$ws.Range("H85").Value = 2042
$ws.Range("I85").Value = 2042
$ws.Range("K85").Value = 6126
$ws.Range("M85").Value = -4722

$ws.Range("H97").Value = 732.5
$ws.Range("J97").Value = 732.5
$ws.Range("L97").Value = 2197.5
$ws.Range("N97").Value = -3189.5

$ws.Range("H125").Value = 1887.7858
$ws.Range("I125").Value = 1391.25
$ws.Range("J125").Value = 2549.8333
$ws.Range("K125").Value = 12521.25
$ws.Range("L125").Value = 22948.4997
$ws.Range("M125").Value = -10061.25
$ws.Range("N125").Value = -27868.4997

$ws.Range("H137").Value = 37324.57
$ws.Range("I137").Value = 1417.75
$ws.Range("K137").Value = 4253.25
$ws.Range("M137").Value = -1703.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23368.32
$ws.Range("I32").Value = 23604.299
$ws.Range("J32").Value = 19671.334
$ws.Range("K32").Value = 23604.299
$ws.Range("L32").Value = 19671.334
$ws.Range("M32").Value = -23317.299
$ws.Range("N32").Value = -20245.334

$ws.Range("H97").Value = 768.4211
$ws.Range("I97").Value = 761.1111
$ws.Range("K97").Value = 761.1111
$ws.Range("M97").Value = -265.1111

$ws.Range("H102").Value = 1384.9166
$ws.Range("I102").Value = 1241.125
$ws.Range("J102").Value = 1672.5
$ws.Range("K102").Value = 1241.125
$ws.Range("L102").Value = 1672.5
$ws.Range("M102").Value = 380.875
$ws.Range("N102").Value = -4916.5

$ws.Range("H132").Value = 13413.581
$ws.Range("I132").Value = 1727.25
$ws.Range("J132").Value = 47410.184
$ws.Range("K132").Value = 5181.75
$ws.Range("L132").Value = 142230.552
$ws.Range("M132").Value = -2651.75
$ws.Range("N132").Value = -147290.552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2294.75
$ws.Range("I20").Value = 1993
$ws.Range("J20").Value = 3200
$ws.Range("K20").Value = 1993
$ws.Range("L20").Value = 3200
$ws.Range("M20").Value = -1746
$ws.Range("N20").Value = -3694

$ws.Range("H94").Value = 1390.1464
$ws.Range("I94").Value = 1021.125
$ws.Range("J94").Value = 2702.2222
$ws.Range("K94").Value = 1021.125
$ws.Range("L94").Value = 2702.2222
$ws.Range("M94").Value = -570.125
$ws.Range("N94").Value = -3604.2222

$ws.Range("H99").Value = 1511.2963
$ws.Range("I99").Value = 1264.7727
$ws.Range("K99").Value = 1264.7727
$ws.Range("M99").Value = 233.2273

$ws.Range("H105").Value = 2384296
$ws.Range("I105").Value = 3734.5833
$ws.Range("J105").Value = 5558378
$ws.Range("K105").Value = 3734.5833
$ws.Range("L105").Value = 5558378
$ws.Range("M105").Value = -1987.5833
$ws.Range("N105").Value = -5561872

$ws.Range("H134").Value = 30350.611
$ws.Range("I134").Value = 35039.418
$ws.Range("J134").Value = 1280
$ws.Range("K134").Value = 105118.254
$ws.Range("L134").Value = 3840
$ws.Range("M134").Value = -102583.254
$ws.Range("N134").Value = -8910

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1829
$ws.Range("I16").Value = 2250
$ws.Range("K16").Value = 2250
$ws.Range("M16").Value = -1963

$ws.Range("H31").Value = 3025.8572
$ws.Range("I31").Value = 1695.7778
$ws.Range("K31").Value = 1695.7778
$ws.Range("M31").Value = -1400.7778

$ws.Range("H34").Value = 3025.8572
$ws.Range("I34").Value = 1695.7778
$ws.Range("K34").Value = 1695.7778
$ws.Range("M34").Value = -1493.7778

$ws.Range("H105").Value = 8334549.5
$ws.Range("J105").Value = 1484.625
$ws.Range("L105").Value = 1484.625
$ws.Range("N105").Value = -4978.625

$ws.Range("H113").Value = 1829
$ws.Range("I113").Value = 2250
$ws.Range("K113").Value = 2250
$ws.Range("M113").Value = -80

$ws.Range("H134").Value = 1168.3043
$ws.Range("I134").Value = 938.2308
$ws.Range("J134").Value = 1467.4
$ws.Range("K134").Value = 2814.6924
$ws.Range("L134").Value = 4402.200000000001
$ws.Range("M134").Value = -279.6923999999999
$ws.Range("N134").Value = -9472.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 266.33334
$ws.Range("I33").Value = 249.5
$ws.Range("K33").Value = 1497
$ws.Range("M33").Value = -1214

$ws.Range("H36").Value = 2619.8572
$ws.Range("I36").Value = 2146.6
$ws.Range("J36").Value = 3803
$ws.Range("K36").Value = 6439.799999999999
$ws.Range("L36").Value = 11409
$ws.Range("M36").Value = -6270.799999999999
$ws.Range("N36").Value = -11747

$ws.Range("H97").Value = 1051.9375
$ws.Range("J97").Value = 1350.125
$ws.Range("L97").Value = 4050.375
$ws.Range("N97").Value = -5042.375

$ws.Range("H131").Value = 765.4545000000001
$ws.Range("J131").Value = 783.163
$ws.Range("L131").Value = 2349.489
$ws.Range("N131").Value = -12429.489

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1536.1538
$ws.Range("I97").Value = 933.1111
$ws.Range("J97").Value = 2893
$ws.Range("K97").Value = 933.1111
$ws.Range("L97").Value = 2893
$ws.Range("M97").Value = -437.1111
$ws.Range("N97").Value = -3885

$ws.Range("H102").Value = 1885.3704
$ws.Range("I102").Value = 1891.875
$ws.Range("K102").Value = 1891.875
$ws.Range("M102").Value = -269.875

$ws.Range("H107").Value = 514.4666999999999
$ws.Range("I107").Value = 224.45454
$ws.Range("J107").Value = 1312
$ws.Range("K107").Value = 224.45454
$ws.Range("L107").Value = 1312
$ws.Range("M107").Value = 1695.54546
$ws.Range("N107").Value = -5152

$ws.Range("H132").Value = 68135.48
$ws.Range("I132").Value = 69758.53
$ws.Range("J132").Value = 65092.25
$ws.Range("K132").Value = 209275.59
$ws.Range("L132").Value = 195276.75
$ws.Range("M132").Value = -206745.59
$ws.Range("N132").Value = -200336.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6300.4
$ws.Range("J22").Value = 6967
$ws.Range("L22").Value = 6967
$ws.Range("N22").Value = -7557

$ws.Range("H27").Value = 6300.4
$ws.Range("J27").Value = 6967
$ws.Range("L27").Value = 6967
$ws.Range("N27").Value = -7181

$ws.Range("H40").Value = 3603.375
$ws.Range("I40").Value = 2521.2856
$ws.Range("K40").Value = 2521.2856
$ws.Range("M40").Value = -2385.2856

$ws.Range("H93").Value = 2465.3333
$ws.Range("I93").Value = 2648.5
$ws.Range("K93").Value = 2648.5
$ws.Range("M93").Value = -1400.5

$ws.Range("H100").Value = 2378.4285
$ws.Range("I100").Value = 2112.5
$ws.Range("J100").Value = 2733
$ws.Range("K100").Value = 2112.5
$ws.Range("L100").Value = 2733
$ws.Range("M100").Value = -1571.5
$ws.Range("N100").Value = -3815

$ws.Range("H136").Value = 1798
$ws.Range("I136").Value = 1629.1428
$ws.Range("K136").Value = 4887.428400000001
$ws.Range("M136").Value = -2337.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H126").Value = 1287.88
$ws.Range("I126").Value = 883.0833
$ws.Range("J126").Value = 1661.5385
$ws.Range("K126").Value = 2649.2499
$ws.Range("L126").Value = 4984.6155
$ws.Range("M126").Value = -179.2498999999998
$ws.Range("N126").Value = -9924.6155
